$d = $word.ActiveDocument

# Locate the target paragraph (the "No Scrum..." paragraph, currently split
# across two runs around a relocated _GoBack bookmark) using Find, then grab
# its enclosing paragraph so we replace the whole thing (bookmark included).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("No Scrum o dono da empresa", $true, $false,
                                    $false, $false, $false, $true, 1, $false,
                                    "", 0)

$target = $null
if ($found) {
    $target = $searchRange.Paragraphs(1)
} else {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text -like "No Scrum o dono da empresa*") {
            $target = $para
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$full = $d.Range($target.Range.Start, $target.Range.End)

$wordml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>No Scrum o dono da empresa solicitando o software seria o Product Owner, figura essa que tem a função de ditar a lista de afazeres, e a ordem de preferência, mas já que a lista está incompleta este método não seria eficiente.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p/>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Equipe</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Líder</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Quanto à criação da equipe de desenvolvimento do software, precisamos primeiramente de um líder da equipe, responsável por observar a equipe e orientá-la.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>O foco é montar a equipe, estabelecer objetivos, papéis e regras; fazer o time interagir e fazer as pessoas “vestirem a camisa</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">”, </w:t>
  </w:r>
  <w:r>
    <w:t>manter a visão crítica do time</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>e utilizar</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> procedimentos estruturados para se comunicar, resolver conflitos, alocar recursos e relacionar-se com a organização</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Funcionários</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">É necessário pessoal para separar os requisitos do software, coletar dados e verificar processos que necessitam de evolução, trabalhando sobre as ordens do líder, os funcionários </w:t>
  </w:r>
  <w:r>
    <w:t>têm</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve"> a função de fazer o software, todos papéis desde levantamento de requisitos á programação devem ser definidos pelo líder da equipe e atribuídos aos funcionários.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p/>
'@

$full.InsertXML($wordml)
